$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -7
$ws.Range("F7").Value = -7
$ws.Range("F9").Value = -6
$ws.Range("F14").Value = -9
$ws.Range("F15").Value = -5
$ws.Range("F19").Value = -4
$ws.Range("F20").Value = -4
$ws.Range("F21").Value = 5
